$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per diff. Cells whose new value is a plain number
# (and would otherwise be auto-converted to a numeric type by Excel) are
# first forced to text format so they keep their exact original text
# representation (matching the source workbook, which stores these as
# literal text/inline strings, e.g. keeping trailing zeros like "1.20").

$ws.Range("D2").Value = "54.260.89"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.265.45"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "496.39"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.94"
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0953"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.338"
$ws.Range("E11").Value = "  +4.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.78"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("D13").Value = "2.665.90"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.88"
$ws.Range("E14").Value = "  +5.15%  "
$ws.Range("D15").Value = "54.244.17"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000130"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "2.254.17"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.25"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.14"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "302.75"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.33"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "60.98"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.997"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  +3.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.69"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0692"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.60"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.95"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.77"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.938"
$ws.Range("E35").Value = "  +8.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.20"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.40"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "124.96"
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0493"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0894"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.547"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "241.15"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0205"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.10"
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("E51").Value = "  -1.05%  "
